$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "[18:30:39]"
$ws.Range("B2").Value = "Свойство и признаки ромба"
$ws.Range("D2").Value = "asdf"
$ws.Range("E2").Value = 1381120546
$ws.Range("F2").Value = "новая команда (2)"
$ws.Range("G2").Value = -300
